$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

foreach ($ws in @($zhcn, $dede)) {
    # Status flips from "Ready for handoff" to "Handoff transform failed"
    $ws.Range("B2").Value = "Handoff transform failed"

    # The handoff attempt produced no handoff file, so the "Latest Handoff
    # File" cell (and its hyperlink to the .xlf) is removed entirely.
    # Deleting a range's Hyperlinks collection clears every hyperlink on the
    # sheet in this engine, so clear first and re-create the two that must
    # survive (the source-file link and the .localization-config link).
    $ws.Range("C2").Hyperlinks.Delete()
    $ws.Range("C2").Clear()

    # Latest Handoff Datetime / Latest Handback DateTime reset to the epoch
    $ws.Range("D2").Value = "0001-01-01 00:00:00"
    $ws.Range("G2").Value = "0001-01-01 00:00:00"

    # Handoff Reason becomes "Ignored", matching row 3's existing value
    $ws.Range("H2").Value = "Ignored"
}

# Re-add the hyperlinks that must remain, pointing at their original targets
$zhcnMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/af9a1137614f95b5bbe2378eaf91b7903ed055cd/e2e/987b8709-0d39-4f67-9612-2a48fc867e1e.md"
$zhcnCfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/af9a1137614f95b5bbe2378eaf91b7903ed055cd/.localization-config"
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $zhcnMdUrl, "", "", "987b8709-0d39-4f67-9612-2a48fc867e1e.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $zhcnCfgUrl, "", "", ".localization-config")

$dedeMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/af9a1137614f95b5bbe2378eaf91b7903ed055cd/e2e/987b8709-0d39-4f67-9612-2a48fc867e1e.md"
$dedeCfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/af9a1137614f95b5bbe2378eaf91b7903ed055cd/.localization-config"
$dede.Hyperlinks.Add($dede.Range("A2"), $dedeMdUrl, "", "", "987b8709-0d39-4f67-9612-2a48fc867e1e.md")
$dede.Hyperlinks.Add($dede.Range("A3"), $dedeCfgUrl, "", "", ".localization-config")
